$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.146.27"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "1.654.70"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.48"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5089"
$ws.Range("E6").Value = "  -4.14%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2580"
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06405"
$ws.Range("E9").Value = "  -4.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.91"
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07785"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").Value = "1.655.36"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.280"
$ws.Range("E13").Value = "  -5.08%  "
$ws.Range("D14").Value = "1.882.62"
$ws.Range("E14").Value = "  -3.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5512"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "0.0₅8009"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.94"
$ws.Range("E17").Value = "  -6.21%  "
$ws.Range("D18").Value = "26.169.83"
$ws.Range("E18").Value = "  -4.29%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "209.22"
$ws.Range("E20").Value = "  -6.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.405"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.05"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.021"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.68"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.736"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1176"
$ws.Range("E27").Value = "  -2.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.974"
$ws.Range("E28").Value = "  -3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.80"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05110"
$ws.Range("E30").Value = "  -4.76%  "
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.346"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.218"
$ws.Range("E33").Value = "  -6.36%  "
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.752"
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9273"
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.364"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").Value = "1.161.57"
$ws.Range("E39").Value = "  +6.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01588"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.554"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.640"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.42"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "1.793.05"
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").Value = "0.0₈116"
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4549"
$ws.Range("E48").Value = "  +0.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.72"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.844"
$ws.Range("E51").Value = "  -3.30%  "
